# Updated price / volume(1h) figures and a row swap (VeChain <-> Bittensor)
# as produced by the scheduled "Updated cryptos list" GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "64.169.27"
    "E2" = "  -0.40%  "
    "D3" = "3.132.19"
    "E3" = "  -1.38%  "
    "E4" = "  -0.02%  "
    "D5" = "569.70"
    "E5" = "  -0.05%  "
    "D6" = "161.08"
    "E6" = "  -4.71%  "
    "E7" = "  -0.16%  "
    "D8" = "0.564"
    "E8" = "  -7.39%  "
    "E9" = "  -3.87%  "
    "D10" = "6.57"
    "E10" = "  -3.13%  "
    "D11" = "0.378"
    "E11" = "  -1.98%  "
    "D12" = "3.674.41"
    "E12" = "  -1.47%  "
    "E13" = "  -0.95%  "
    "D14" = "64.217.16"
    "E14" = "  -0.39%  "
    "D15" = "24.82"
    "E15" = "  -2.18%  "
    "D16" = "3.134.52"
    "E16" = "  -1.50%  "
    "E17" = "  -3.34%  "
    "D18" = "400.17"
    "E18" = "  -4.78%  "
    "D19" = "12.47"
    "E19" = "  -3.06%  "
    "D20" = "5.20"
    "E20" = "  -3.16%  "
    "D21" = "7.08"
    "E21" = "  +0.33%  "
    "D22" = "5.86"
    "E22" = "  +3.46%  "
    "D23" = "0.999"
    "E23" = "  -0.13%  "
    "D24" = "67.65"
    "E24" = "  -3.71%  "
    "D25" = "0.480"
    "E25" = "  -1.62%  "
    "E26" = "  -5.13%  "
    "D27" = "0.0₃0999"
    "E27" = "  -5.53%  "
    "D28" = "8.77"
    "E28" = "  +0.04%  "
    "D29" = "0.997"
    "E29" = "  -0.27%  "
    "E30" = "  +0.08%  "
    "E31" = "  -1.79%  "
    "D32" = "21.08"
    "E32" = "  -3.15%  "
    "D33" = "158.42"
    "E33" = "  +0.58%  "
    "E34" = "  -1.81%  "
    "D35" = "4.77"
    "E35" = "  -5.23%  "
    "E36" = "  -3.04%  "
    "E37" = "  -3.01%  "
    "D38" = "2.657.48"
    "E38" = "  -2.31%  "
    "E39" = "  -2.71%  "
    "D40" = "23.35"
    "E40" = "  -3.99%  "
    "D41" = "4.05"
    "E41" = "  -2.69%  "
    "D42" = "38.33"
    "E42" = "  -2.05%  "
    "D43" = "0.685"
    "E43" = "  -3.53%  "
    "E44" = "  -2.13%  "
    "D45" = "5.43"
    "E45" = "  -3.08%  "
    "B46" = "Bittensor"
    "C46" = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
    "D46" = "288.03"
    "E46" = "  -1.51%  "
    "B47" = "VeChain"
    "C47" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "D47" = "0.0253"
    "E47" = "  -2.93%  "
    "D48" = "20.93"
    "E48" = "  -3.51%  "
    "D49" = "0.996"
    "E49" = "  -0.37%  "
    "E50" = "  -1.97%  "
    "D51" = "10.45"
    "E51" = "  -0.08%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Leading apostrophe forces text entry so Excel doesn't
    # reinterpret numeric-looking values (e.g. "569.70", "5.20")
    # and strip significant trailing zeros.
    $cell.Value = "'" + $updates[$addr]
    # Reset style to Normal so the quote-prefix flag added above
    # does not leave a stray cell style behind.
    $cell.Style = "Normal"
}
